$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value renders as plain numeric text need to be forced back to
# text (matching the source inlineStr cells) since Excel auto-converts numeric-
# looking strings assigned via .Value into real numbers.
function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

$ws.Range("D2").Value = "30.167.90"
$ws.Range("E2").Value = "  -4.21%  "
$ws.Range("D3").Value = "1.917.88"
$ws.Range("E3").Value = "  -3.37%  "
$ws.Range("E4").Value = "  +0.90%  "
Set-TextValue $ws.Range("D5") "245.55"
$ws.Range("E5").Value = "  -2.59%  "
Set-TextValue $ws.Range("D6") "0.6967"
$ws.Range("E6").Value = "  -15.32%  "
Set-TextValue $ws.Range("D7") "1.003"
$ws.Range("E7").Value = "  +0.97%  "
Set-TextValue $ws.Range("D8") "0.3226"
$ws.Range("E8").Value = "  -6.21%  "
Set-TextValue $ws.Range("D9") "26.13"
$ws.Range("E9").Value = "  +0.98%  "
Set-TextValue $ws.Range("D10") "0.06796"
$ws.Range("E10").Value = "  -2.36%  "
Set-TextValue $ws.Range("D11") "0.7895"
$ws.Range("E11").Value = "  -5.82%  "
Set-TextValue $ws.Range("D12") "0.07929"
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("D13").Value = "1.925.73"
$ws.Range("E13").Value = "  -2.74%  "
Set-TextValue $ws.Range("D14") "5.353"
$ws.Range("E14").Value = "  -2.56%  "
Set-TextValue $ws.Range("D15") "93.78"
$ws.Range("E15").Value = "  -8.34%  "
Set-TextValue $ws.Range("D16") "14.36"
$ws.Range("E16").Value = "  +2.59%  "
Set-TextValue $ws.Range("D17") "259.01"
$ws.Range("E17").Value = "  -5.90%  "
$ws.Range("D18").Value = "30.192.21"
$ws.Range("E18").Value = "  -4.05%  "
Set-TextValue $ws.Range("D19") "5.796"
$ws.Range("E19").Value = "  +1.94%  "
Set-TextValue $ws.Range("D20") "0.000007805"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("D21").Value = "2.179.28"
$ws.Range("E21").Value = "  -3.00%  "
$ws.Range("E22").Value = "  +0.82%  "
Set-TextValue $ws.Range("D23") "1.004"
$ws.Range("E23").Value = "  +0.98%  "
Set-TextValue $ws.Range("D24") "6.786"
$ws.Range("E24").Value = "  -1.37%  "
Set-TextValue $ws.Range("D25") "9.565"
$ws.Range("E25").Value = "  -0.93%  "
Set-TextValue $ws.Range("D26") "159.75"
$ws.Range("E26").Value = "  -3.53%  "
Set-TextValue $ws.Range("D27") "18.64"
$ws.Range("E27").Value = "  -5.61%  "
Set-TextValue $ws.Range("D28") "0.1305"
$ws.Range("E28").Value = "  -19.26%  "
Set-TextValue $ws.Range("D29") "2.210"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D30") "1.348"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.547"
$ws.Range("E31").Value = "  -0.40%  "
Set-TextValue $ws.Range("D32") "4.385"
$ws.Range("E32").Value = "  -3.94%  "
Set-TextValue $ws.Range("D33") "4.155"
$ws.Range("E33").Value = "  -3.91%  "
Set-TextValue $ws.Range("D34") "0.05023"
$ws.Range("E34").Value = "  -3.39%  "
Set-TextValue $ws.Range("D35") "1.183"
$ws.Range("E35").Value = "  -2.99%  "
Set-TextValue $ws.Range("D36") "0.7411"
$ws.Range("E36").Value = "  -0.54%  "
Set-TextValue $ws.Range("D37") "2.713"
$ws.Range("E37").Value = "  -1.64%  "
Set-TextValue $ws.Range("D38") "0.01914"
$ws.Range("E38").Value = "  -3.40%  "
Set-TextValue $ws.Range("D39") "2.788"
$ws.Range("E39").Value = "  -4.18%  "
Set-TextValue $ws.Range("D40") "79.62"
$ws.Range("E40").Value = "  +1.44%  "
Set-TextValue $ws.Range("D41") "6.474"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("E42").Value = "  -3.24%  "
Set-TextValue $ws.Range("D43") "0.4370"
$ws.Range("E43").Value = "  -6.02%  "
Set-TextValue $ws.Range("D44") "1.002"
$ws.Range("E44").Value = "  +0.95%  "
Set-TextValue $ws.Range("D45") "0.8336"
$ws.Range("E45").Value = "  -2.24%  "
Set-TextValue $ws.Range("D46") "101.55"
$ws.Range("E46").Value = "  -3.70%  "
Set-TextValue $ws.Range("D47") "9.628"
$ws.Range("E47").Value = "  -3.02%  "
Set-TextValue $ws.Range("D48") "7.144"
$ws.Range("E48").Value = "  -5.18%  "
Set-TextValue $ws.Range("D49") "35.68"
$ws.Range("E49").Value = "  -2.05%  "
Set-TextValue $ws.Range("D50") "0.05925"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D51") "1.467"
$ws.Range("E51").Value = "  +1.67%  "
